# This script applies a re-crawl refresh to the product export sheet:
#  - the timestamp column (O) is bumped to the new crawl time for every
#    data row
#  - three pairs of product rows that the crawler emitted in a different
#    order this run are swapped back into their new positions
#    (6<->7, 9<->10, 21<->22)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2022-08-06 07:01:09"
$newTimestamp = "2022-08-06 20:57:37"

# Columns that can contain values which *look* numeric (ids, prices) but
# must stay stored as text, matching the original export format.
$textNumericCols = @("A", "H", "K")

function Swap-Rows($rowA, $rowB) {
    $lastCol = 14   # column N; columns A..N hold the product data
    $valsA = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0)
    $valsB = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0)

    for ($c = 1; $c -le $lastCol; $c++) {
        $valsA[$c - 1] = $ws.Cells.Item($rowA, $c).Value2
        $valsB[$c - 1] = $ws.Cells.Item($rowB, $c).Value2
    }

    for ($c = 1; $c -le $lastCol; $c++) {
        $colLetter = $ws.Cells.Item(1, $c).Address($false, $false) -replace '[0-9]', ''

        $cellA = $ws.Cells.Item($rowA, $c)
        $cellB = $ws.Cells.Item($rowB, $c)

        $newValForA = $valsB[$c - 1]
        $newValForB = $valsA[$c - 1]

        if ($textNumericCols -contains $colLetter) {
            $cellA.NumberFormat = "@"
            $cellB.NumberFormat = "@"
        }

        if ($newValForA -ne $null) { $cellA.Value2 = $newValForA }
        if ($newValForB -ne $null) { $cellB.Value2 = $newValForB }
    }
}

Swap-Rows 6 7
Swap-Rows 9 10
Swap-Rows 21 22

# Refresh the crawl timestamp on every data row (row 1 is the header).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 15)   # column O
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value2 = $newTimestamp
    }
}
